# Auto-generated Excel COM-interop script
# Applies "Add data for 2024-10-18" updates to 2024 (column K) values
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6243
$ws.Range("K3").Value = 6449
$ws.Range("K4").Value = 1346
$ws.Range("K5").Value = 460
$ws.Range("K6").Value = 7101
$ws.Range("K7").Value = 21599

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K6").Value = 116
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 387
$ws.Range("K3").Value = 434
$ws.Range("K6").Value = 480
$ws.Range("K7").Value = 1421

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 469

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 243
$ws.Range("K3").Value = 337
$ws.Range("K6").Value = 287
$ws.Range("K7").Value = 939

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 126
$ws.Range("K7").Value = 355

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 206
$ws.Range("K3").Value = 244
$ws.Range("K6").Value = 214
$ws.Range("K7").Value = 730

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 127
$ws.Range("K6").Value = 181
$ws.Range("K7").Value = 509

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 147
$ws.Range("K7").Value = 354

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 149
$ws.Range("K7").Value = 634
$ws.Range("K8").Value = 1421
$ws.Range("K9").Value = 94
$ws.Range("K10").Value = 124
$ws.Range("K11").Value = 401
$ws.Range("K14").Value = 109
$ws.Range("K15").Value = 220
$ws.Range("K19").Value = 630
$ws.Range("K25").Value = 103
$ws.Range("K27").Value = 205
$ws.Range("K29").Value = 1172
$ws.Range("K31").Value = 239
$ws.Range("K33").Value = 939
$ws.Range("K36").Value = 277
$ws.Range("K37").Value = 730
$ws.Range("K41").Value = 150
$ws.Range("K42").Value = 798
$ws.Range("K45").Value = 30
$ws.Range("K46").Value = 43
$ws.Range("K50").Value = 102
$ws.Range("K52").Value = 567
$ws.Range("K53").Value = 275
$ws.Range("K54").Value = 424
$ws.Range("K55").Value = 237
$ws.Range("K63").Value = 64
$ws.Range("K64").Value = 136
$ws.Range("K65").Value = 509
$ws.Range("K69").Value = 48
$ws.Range("K70").Value = 37
$ws.Range("K72").Value = 109
$ws.Range("K73").Value = 194
$ws.Range("K76").Value = 295
$ws.Range("K79").Value = 545
$ws.Range("K83").Value = 469
$ws.Range("K84").Value = 173
$ws.Range("K85").Value = 1006
$ws.Range("K89").Value = 317
$ws.Range("K91").Value = 250
$ws.Range("K95").Value = 355
$ws.Range("K99").Value = 354
$ws.Range("K101").Value = 21599

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 77
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 67
$ws.Range("K7").Value = 173

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 102
$ws.Range("K6").Value = 231
$ws.Range("K7").Value = 424

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 334
$ws.Range("K5").Value = 29
$ws.Range("K7").Value = 1172

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 205
$ws.Range("K7").Value = 630

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 295

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 56
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 219
$ws.Range("K7").Value = 798

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 66
$ws.Range("K4").Value = 11
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 237

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 120
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 250

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 176
$ws.Range("K6").Value = 138
$ws.Range("K7").Value = 545

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 27
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 109
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 207
$ws.Range("K6").Value = 173
$ws.Range("K7").Value = 634

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 79
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K6").Value = 130
$ws.Range("K7").Value = 401

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 317

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K6").Value = 250
$ws.Range("K7").Value = 1006

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 202
$ws.Range("K7").Value = 567

Write-Host "Applied 2024-10-18 update to $($wb.Worksheets.Count) worksheets (144 cells)."
